$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R1").Value = "DINAMICAS"
$ws.Range("R7").Value = "ESTATICAS"
$ws.Range("R13").Value = "TRANSICAO"
$ws.Range("H7").Value = "DONT NEED"

$ws.Range("H8").Select()
